$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Sheet "AMSIN" (sheet1): row 22 gets the same "Normal" formatting
#    that every other data row (e.g. row 21) already carries, and the
#    run-time value in B22 gets corrected to its precise serial value.
# -----------------------------------------------------------------
$amsin = $wb.Worksheets.Item("AMSIN")

foreach ($col in @("A", "C", "D", "E", "F", "G")) {
    $amsin.Range($col + "22").Style = "Normal"
}

$amsin.Range("B22").Value = 44810.94245770833

# -----------------------------------------------------------------
# 2) Sheet "AMS" (sheet3): append a new row 15 with the latest
#    registration run ("ocr166"), growing the sheet's used range
#    from A1:G14 to A1:G15.
# -----------------------------------------------------------------
$ams = $wb.Worksheets.Item("AMS")

# Column A holds a date-like label ("2022-09-08") that must stay a
# literal text value instead of being auto-converted to a date serial.
$ams.Range("A15").NumberFormat = "@"
$ams.Range("A15").Value = "2022-09-08"

# Column B is the precise run timestamp (serial date/time number),
# formatted the same way the other sheets' "Run Time" column is.
$ams.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ams.Range("B15").Value = 44812.54125285153

$ams.Range("C15").Value = "ocr166"
$ams.Range("D15").Value = 42
$ams.Range("E15").Value = 42
$ams.Range("F15").Value = 0
$ams.Range("G15").Value = 1.09
